$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 2: NEI-CEDA CI
$ws1.Range("D2").Value = 3670
$ws1.Range("E2").Value = 945

# Row 3: BRVM - SERVICES PUBLICS
$ws1.Range("D3").Value = 3308.27
$ws1.Range("E3").Value = 105.42

# Row 4: BRVM - AUTRES SECTEURS
$ws1.Range("D4").Value = 2419.67
$ws1.Range("E4").Value = 623.11

# Row 5: BRVM - DISTRIBUTION
$ws1.Range("D5").Value = 1987.5
$ws1.Range("E5").Value = 486.27

# Row 6: BRVM - TRANSPORT
$ws1.Range("D6").Value = 1425.8
$ws1.Range("E6").Value = 351.25

# Row 7: BRVM - AGRICULTURE
$ws1.Range("D7").Value = 1328.71
$ws1.Range("E7").Value = 325.09

# Row 8: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Range("D8").Value = 680.17
$ws1.Range("E8").Value = 167.85

# Row 9: BRVM - FINANCES
$ws1.Range("D9").Value = 584.33
$ws1.Range("E9").Value = 145.2

# Row 10: BRVM - SERVICES FINANCIERS
$ws1.Range("D10").Value = 574.28
$ws1.Range("E10").Value = 142.7

# Row 11: BRVM-PRESTIGE
$ws1.Range("D11").Value = 561.34
$ws1.Range("E11").Value = 139.52

# Row 12: BRVM - INDUSTRIELS
$ws1.Range("D12").Value = 514.8099999999999
$ws1.Range("E12").Value = 130.26

# Row 13: BRVM - INDUSTRIE  (**)
$ws1.Range("A13").Value = 'BRVM - INDUSTRIE  (**)'
$ws1.Range("C13").Value = 2
$ws1.Range("D13").Value = 514.38
$ws1.Range("E13").Value = 257.06

# Row 14: BRVM - ENERGIE
$ws1.Range("A14").Value = 'BRVM - ENERGIE'
$ws1.Range("D14").Value = 444.42
$ws1.Range("E14").Value = 109.09

# Row 15: BRVM-PRINCIPAL  (**)
$ws1.Range("A15").Value = 'BRVM-PRINCIPAL  (**)'
$ws1.Range("C15").Value = 2
$ws1.Range("D15").Value = 432.54
$ws1.Range("E15").Value = 215.71

# Row 16: BRVM - CONSOMMATION DE BASE  (**)
$ws1.Range("A16").Value = 'BRVM - CONSOMMATION DE BASE  (**)'
$ws1.Range("C16").Value = 2
$ws1.Range("D16").Value = 427.68
$ws1.Range("E16").Value = 213.14

# Row 17: BRVM - TELECOMMUNICATIONS
$ws1.Range("A17").Value = 'BRVM - TELECOMMUNICATIONS'
$ws1.Range("C17").Value = 4
$ws1.Range("D17").Value = 373.78
$ws1.Range("E17").Value = 93.48

# Row 18: BRVM - INDUSTRIE
$ws1.Range("A18").Value = 'BRVM - INDUSTRIE'
$ws1.Range("D18").Value = 266.4
$ws1.Range("E18").Value = 266.4

# Row 19: BRVM - INDUSTRIE    (**)
$ws1.Range("A19").Value = 'BRVM - INDUSTRIE    (**)'
$ws1.Range("D19").Value = 262.27
$ws1.Range("E19").Value = 262.27

# Row 20: BRVM - CONSOMMATION DE BASE
$ws1.Range("A20").Value = 'BRVM - CONSOMMATION DE BASE'
$ws1.Range("D20").Value = 222.06
$ws1.Range("E20").Value = 222.06

# Row 21: BRVM-PRINCIPAL
$ws1.Range("A21").Value = 'BRVM-PRINCIPAL'
$ws1.Range("D21").Value = 220.02
$ws1.Range("E21").Value = 220.02

# Row 22: BRVM-PRINCIPAL     (**)
$ws1.Range("A22").Value = 'BRVM-PRINCIPAL     (**)'
$ws1.Range("D22").Value = 219.45
$ws1.Range("E22").Value = 219.45

# Row 23: BRVM - CONSOMMATION DE BASE   (**)
$ws1.Range("A23").Value = 'BRVM - CONSOMMATION DE BASE   (**)'
$ws1.Range("D23").Value = 218.71
$ws1.Range("E23").Value = 218.71

# Row 24: ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)
$ws1.Range("D24").Value = 29.53
$ws1.Range("E24").Value = 7.25

# Row 25: SETAO CI (STAC)
$ws1.Range("B25").Value = 2
$ws1.Range("D25").Value = 14.8
$ws1.Range("E25").Value = 7.36

# Row 26: NEI-CEDA CI (NEIC)
$ws1.Range("A26").Value = 'NEI-CEDA CI (NEIC)'
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = 5.5
$ws1.Range("E26").Value = 4.76
$ws1.Range("G26").Value = '👀 À surveiller'

# Row 27: BANK OF AFRICA NG (BOAN)
$ws1.Range("A27").Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = 4.42
$ws1.Range("E27").Value = 5.77
$ws1.Range("G27").Value = '👀 À surveiller'

# Row 28: SMB CI (SMBC)
$ws1.Range("A28").Value = 'SMB CI (SMBC)'
$ws1.Range("D28").Value = 3.19
$ws1.Range("E28").Value = 3.19

# Row 29: ECOBANK COTE D''IVOIRE (ECOC)
$ws1.Range("A29").Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Range("D29").Value = 2.89
$ws1.Range("E29").Value = 2.89

# Row 30: AFRICA GLOBAL LOGISTICS CI (SDSC)
$ws1.Range("A30").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Range("D30").Value = 2.76
$ws1.Range("E30").Value = 2.76

# Row 31: ORAGROUP TOGO (ORGT)
$ws1.Range("A31").Value = 'ORAGROUP TOGO (ORGT)'
$ws1.Range("C31").Value = 0
$ws1.Range("D31").Value = 2.6
$ws1.Range("E31").Value = 2.6
$ws1.Range("G31").Value = '➖ Neutre'

# Row 32: TOTALENERGIES MARKETING SN (TTLS)
$ws1.Range("A32").Value = 'TOTALENERGIES MARKETING SN (TTLS)'
$ws1.Range("C32").Value = 0
$ws1.Range("D32").Value = 2.41
$ws1.Range("E32").Value = 2.41
$ws1.Range("G32").Value = '➖ Neutre'

# Row 33: ONATEL BF (ONTBF)
$ws1.Range("A33").Value = 'ONATEL BF (ONTBF)'
$ws1.Range("D33").Value = 2.04
$ws1.Range("E33").Value = 2.04

# Row 34: BERNABE CI (BNBC)
$ws1.Range("A34").Value = 'BERNABE CI (BNBC)'
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = 1.83
$ws1.Range("E34").Value = -1.38
$ws1.Range("G34").Value = '👀 À surveiller'

# Row 35: TRACTAFRIC MOTORS CI (PRSC)
$ws1.Range("A35").Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws1.Range("C35").Value = 0
$ws1.Range("D35").Value = 1.6
$ws1.Range("E35").Value = 1.6
$ws1.Range("G35").Value = '➖ Neutre'

# Row 36: SICABLE CI (CABC)
$ws1.Range("A36").Value = 'SICABLE CI (CABC)'
$ws1.Range("B36").Value = 1
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = 1.32
$ws1.Range("E36").Value = 7.5
$ws1.Range("G36").Value = '👀 À surveiller'

# Row 37: TOTAL
$ws1.Range("A37").Value = 'TOTAL'
$ws1.Range("B37").Value = 0
$ws1.Range("C37").Value = 4
$ws1.Range("D37").Value = 0
$ws1.Range("E37").Value = 0
$ws1.Range("G37").Value = '➖ Neutre'

# Row 38: SUCRIVOIRE (SCRC)
$ws1.Range("A38").Value = 'SUCRIVOIRE (SCRC)'
$ws1.Range("B38").Value = 1
$ws1.Range("D38").Value = -0.38
$ws1.Range("E38").Value = -3.14
$ws1.Range("G38").Value = '👀 À surveiller'

# Row 39: BICI CI (BICC)
$ws1.Range("A39").Value = 'BICI CI (BICC)'
$ws1.Range("D39").Value = -1.64
$ws1.Range("E39").Value = -1.64

# Row 40: BANK OF AFRICA CI (BOAC)
$ws1.Range("A40").Value = 'BANK OF AFRICA CI (BOAC)'
$ws1.Range("B40").Value = 0
$ws1.Range("C40").Value = 1
$ws1.Range("D40").Value = -2.3
$ws1.Range("E40").Value = -2.3
$ws1.Range("G40").Value = '➖ Neutre'

# Row 41: CFAO MOTORS CI (CFAC)
$ws1.Range("A41").Value = 'CFAO MOTORS CI (CFAC)'
$ws1.Range("D41").Value = -3.45
$ws1.Range("E41").Value = -3.45

# Row 42: NSIA BANQUE COTE D'IVOIRE (NSBC)

# Row 43: ORANGE COTE D'IVOIRE (ORAC)

# Row 44: CIE CI (CIEC)
$ws1.Range("A44").Value = 'CIE CI (CIEC)'
$ws1.Range("D44").Value = -4.09
$ws1.Range("E44").Value = -4.09

# Row 45: SERVAIR ABIDJAN CI (ABJC)
$ws1.Range("A45").Value = 'SERVAIR ABIDJAN CI (ABJC)'
$ws1.Range("D45").Value = -4.73
$ws1.Range("E45").Value = -4.73

# Row 46: SODE CI (SDCC)
$ws1.Range("A46").Value = 'SODE CI (SDCC)'
$ws1.Range("D46").Value = -4.92
$ws1.Range("E46").Value = -4.92

# Row 47: SAPH CI (SPHC)
$ws1.Range("A47").Value = 'SAPH CI (SPHC)'
$ws1.Range("D47").Value = -5.66
$ws1.Range("E47").Value = -5.66

# Row 48: SICOR CI (SICC)
$ws1.Range("C48").Value = 2
$ws1.Range("D48").Value = -7.28
$ws1.Range("F48").Value = '🟡 Observer'
$ws1.Range("G48").Value = '👀 À surveiller'

# Row 49: VIVO ENERGY CI (SHEC)
$ws1.Range("A49").Value = 'VIVO ENERGY CI (SHEC)'
$ws1.Range("D49").Value = -12.11
$ws1.Range("E49").Value = -6.74

# Row 50: UNILEVER CI (UNLC)
$ws1.Range("A50").Value = 'UNILEVER CI (UNLC)'
$ws1.Range("B50").Value = 0
$ws1.Range("C50").Value = 2
$ws1.Range("D50").Value = -15
$ws1.Range("E50").Value = -7.5
$ws1.Range("F50").Value = '🟡 Observer'
$ws1.Range("G50").Value = '➖ Neutre'
# Top_YTD sheet updates
$ws2.Range("B2").Value = 8574576.5
$ws2.Range("B3").Value = 1070868.25
$ws2.Range("B4").Value = 246701.76
$ws2.Range("B5").Value = 126789.04
$ws2.Range("B6").Value = 43286.66
$ws2.Range("B7").Value = 34777.72
$ws2.Range("B8").Value = 5217.52
$ws2.Range("B9").Value = 3567.02
$ws2.Range("B10").Value = 3419.54
$ws2.Range("B11").Value = 3236.25